$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new (blank) column N, shifting the
# existing "Late" / "Outstanding" / "Disbursement" columns one to the right ---
$ws3 = $wb.Worksheets.Item("Repayment schedule")
$ws3.Columns.Item(14).Insert() | Out-Null
$ws3.Columns.Item(14).ColumnWidth = 9.83

# Make "Repayment schedule" the active / selected tab, with L16 selected
# (matching the cell selection left behind by the edit)
$ws3.Activate() | Out-Null
$ws3.Range("L16").Select() | Out-Null
